# Updates the cryptos price/volume table with refreshed figures (and a
# couple of rows whose coins swapped rank position), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text (preventing Excel from
# re-interpreting dotted numbers like "29.495.52" or "0.9985" as numbers
# or dates), while keeping the cell's format as plain "Normal"/General,
# just like the original inline-string cells.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '29.495.52'
Set-TextValue 'E2' '  -1.03%  '

Set-TextValue 'D3' '1.849.54'
Set-TextValue 'E3' '  -0.61%  '

Set-TextValue 'D4' '0.9985'

Set-TextValue 'D5' '241.82'
Set-TextValue 'E5' '  -0.80%  '

Set-TextValue 'D6' '0.6292'
Set-TextValue 'E6' '  -2.34%  '

Set-TextValue 'D7' '0.9996'
Set-TextValue 'E7' '  -0.03%  '

Set-TextValue 'D8' '0.07541'
Set-TextValue 'E8' '  -0.05%  '

Set-TextValue 'E9' '  -0.40%  '

Set-TextValue 'D10' '24.41'
Set-TextValue 'E10' '  -1.33%  '

Set-TextValue 'E11' '  +0.28%  '

Set-TextValue 'D12' '1.896.68'
Set-TextValue 'E12' '  +1.82%  '

Set-TextValue 'B13' 'Polygon'
Set-TextValue 'C13' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D13' '0.6915'
Set-TextValue 'E13' '  -0.17%  '

Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '5.014'
Set-TextValue 'E14' '  -0.86%  '

Set-TextValue 'D15' '83.71'
Set-TextValue 'E15' '  -0.51%  '

Set-TextValue 'D16' '0.000009799'
Set-TextValue 'E16' '  -0.73%  '

Set-TextValue 'D17' '2.157.75'
Set-TextValue 'E17' '  +1.94%  '

Set-TextValue 'D18' '6.238'
Set-TextValue 'E18' '  +2.01%  '

Set-TextValue 'D19' '29.567.94'
Set-TextValue 'E19' '  -0.80%  '

Set-TextValue 'D20' '234.17'
Set-TextValue 'E20' '  -1.00%  '

Set-TextValue 'D22' '0.9999'
Set-TextValue 'E22' '  +0.00%  '

Set-TextValue 'D23' '7.656'
Set-TextValue 'E23' '  +0.81%  '

Set-TextValue 'D24' '0.9996'
Set-TextValue 'E24' '  -0.06%  '

Set-TextValue 'D25' '154.95'
Set-TextValue 'E25' '  -2.23%  '

Set-TextValue 'E26' '  -2.20%  '

Set-TextValue 'D27' '8.457'
Set-TextValue 'E27' '  -1.38%  '

Set-TextValue 'D28' '17.71'
Set-TextValue 'E28' '  -1.30%  '

Set-TextValue 'D29' '1.478'
Set-TextValue 'E29' '  -0.77%  '

Set-TextValue 'D30' '0.05859'
Set-TextValue 'E30' '  -5.48%  '

Set-TextValue 'D31' '1.252'
Set-TextValue 'E31' '  -2.80%  '

Set-TextValue 'D32' '4.104'
Set-TextValue 'E32' '  -1.34%  '

Set-TextValue 'D33' '4.049'
Set-TextValue 'E33' '  -1.20%  '

Set-TextValue 'D34' '1.882'
Set-TextValue 'E34' '  -0.70%  '

Set-TextValue 'E35' '  -0.40%  '

Set-TextValue 'D36' '0.7213'
Set-TextValue 'E36' '  -1.59%  '

Set-TextValue 'D37' '2.587'
Set-TextValue 'E37' '  -0.88%  '

Set-TextValue 'D38' '1.243.62'
Set-TextValue 'E38' '  +1.99%  '

Set-TextValue 'D39' '2.796'
Set-TextValue 'E39' '  -1.14%  '

Set-TextValue 'D40' '0.01784'
Set-TextValue 'E40' '  -0.47%  '

Set-TextValue 'D41' '0.9051'
Set-TextValue 'E41' '  -1.76%  '

Set-TextValue 'D42' '6.169'
Set-TextValue 'E42' '  -2.52%  '

Set-TextValue 'D43' '2.066.58'
Set-TextValue 'E43' '  +1.59%  '

Set-TextValue 'D44' '0.9993'
Set-TextValue 'E44' '  -0.10%  '

Set-TextValue 'D45' '101.93'
Set-TextValue 'E45' '  -0.13%  '

Set-TextValue 'E46' '  +0.18%  '

Set-TextValue 'D47' '7.376'
Set-TextValue 'E47' '  +9.40%  '

Set-TextValue 'B48' 'TheSandbox'
Set-TextValue 'C48' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D48' '0.4047'
Set-TextValue 'E48' '  -0.65%  '

Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '9.127'
Set-TextValue 'E49' '  -0.85%  '

Set-TextValue 'B50' 'RenderToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '1.708'
Set-TextValue 'E50' '  +2.12%  '

Set-TextValue 'B51' 'BabyDogeCoin'
Set-TextValue 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D51' '0.00000000117'
Set-TextValue 'E51' '  -2.43%  '
